$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Remove the stale selection by reselecting A1 (clears any prior selection ref)
$ws.Range("A1").Select()

# Set column C width (diff adds <col min="3" max="3" width="48" customWidth="1"/>)
# ColumnWidth property adds ~5/6 padding when stored as the OOXML "width" attribute,
# so back that out to land exactly on width="48".
$ws.Columns.Item(3).ColumnWidth = 48 - (5/6)

# copy styles from row 6 (A6:B6) since A7 style matches s="4" (like A4/A5/A6) and B7 matches s="8" (like B5/B6)
$ws.Range("A6:B6").Copy($ws.Range("A7:B7"))
$excel.CutCopyMode = $false

# Add new row of data
$ws.Range("A7").Value = 44499
$ws.Range("B7").Value = 1.45
$ws.Range("C7").Value = "Präsentation, Verschönerungen und README"
